$wb = $excel.ActiveWorkbook

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1104.8235  # H15: 1173.625 -> 1104.8235
$ws.Cells.Item(15, 9).Value = 1104.8235  # I15: 1173.625 -> 1104.8235
$ws.Cells.Item(15, 11).Value = 3314.4705  # K15: 3520.875 -> 3314.4705
$ws.Cells.Item(15, 13).Value = -3145.4705  # M15: -3351.875 -> -3145.4705

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 498.5  # H74: 0 -> 498.5
$ws.Cells.Item(74, 9).Value = 498.5  # I74: 0 -> 498.5
$ws.Cells.Item(74, 11).Value = 498.5  # K74: 0 -> 498.5
$ws.Cells.Item(74, 13).Value = 437.5  # M74: None -> 437.5

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(77, 8).Value = 498.5  # H77: 0 -> 498.5
$ws.Cells.Item(77, 9).Value = 498.5  # I77: 0 -> 498.5
$ws.Cells.Item(77, 11).Value = 2492.5  # K77: 0 -> 2492.5
$ws.Cells.Item(77, 13).Value = 2187.5  # M77: None -> 2187.5

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 25003000  # H113: 36668668 -> 25003000
$ws.Cells.Item(113, 9).Value = 4500  # I113: 5003000 -> 4500
$ws.Cells.Item(113, 10).Value = 50001500  # J113: 100000000 -> 50001500
$ws.Cells.Item(113, 11).Value = 4500  # K113: 5003000 -> 4500
$ws.Cells.Item(113, 12).Value = 50001500  # L113: 100000000 -> 50001500
$ws.Cells.Item(113, 13).Value = -1246  # M113: -4999746 -> -1246
$ws.Cells.Item(113, 14).Value = -50008008  # N113: -100006508 -> -50008008

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 3247.5  # H125: 4165 -> 3247.5
$ws.Cells.Item(125, 9).Value = 2330.3333  # I125: 5997 -> 2330.3333
$ws.Cells.Item(125, 10).Value = 5999  # J125: 3249 -> 5999
$ws.Cells.Item(125, 11).Value = 20972.9997  # K125: 53973 -> 20972.9997
$ws.Cells.Item(125, 12).Value = 53991  # L125: 29241 -> 53991
$ws.Cells.Item(125, 13).Value = -18512.9997  # M125: -51513 -> -18512.9997
$ws.Cells.Item(125, 14).Value = -58911  # N125: -34161 -> -58911

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1125.75  # H132: 969.5 -> 1125.75
$ws.Cells.Item(132, 9).Value = 1125.75  # I132: 969.5 -> 1125.75
$ws.Cells.Item(132, 11).Value = 3377.25  # K132: 2908.5 -> 3377.25
$ws.Cells.Item(132, 13).Value = -847.25  # M132: -378.5 -> -847.25

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 1861.3334  # H135: 1815.2 -> 1861.3334
$ws.Cells.Item(135, 10).Value = 0  # J135: 1400 -> 0
$ws.Cells.Item(135, 12).Value = 0  # L135: 12600 -> 0
$ws.Cells.Item(135, 14).ClearContents()  # N135: -17670 -> (removed)

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 4206.7144  # H138: 4117.722 -> 4206.7144
$ws.Cells.Item(138, 9).Value = 2641.7778  # I138: 2477.6 -> 2641.7778
$ws.Cells.Item(138, 10).Value = 4748.423  # J138: 4748.5386 -> 4748.423
$ws.Cells.Item(138, 11).Value = 7925.3334  # K138: 7432.799999999999 -> 7925.3334
$ws.Cells.Item(138, 12).Value = 14245.269  # L138: 14245.6158 -> 14245.269
$ws.Cells.Item(138, 13).Value = -2785.3334  # M138: -2292.799999999999 -> -2785.3334
$ws.Cells.Item(138, 14).Value = -24525.269  # N138: -24525.6158 -> -24525.269

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4802.143  # H32: 4802.2856 -> 4802.143
$ws.Cells.Item(32, 9).Value = 3825.423  # I32: 3825.577 -> 3825.423
$ws.Cells.Item(32, 11).Value = 3825.423  # K32: 3825.577 -> 3825.423
$ws.Cells.Item(32, 13).Value = -3538.423  # M32: -3538.577 -> -3538.423

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2749.2222  # H61: 2578.3845 -> 2749.2222
$ws.Cells.Item(61, 9).Value = 2593.125  # I61: 2392.818 -> 2593.125
$ws.Cells.Item(61, 10).Value = 3998  # J61: 3599 -> 3998
$ws.Cells.Item(61, 11).Value = 2593.125  # K61: 2392.818 -> 2593.125
$ws.Cells.Item(61, 12).Value = 3998  # L61: 3599 -> 3998
$ws.Cells.Item(61, 13).Value = -2381.125  # M61: -2180.818 -> -2381.125
$ws.Cells.Item(61, 14).Value = -4422  # N61: -4023 -> -4422

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 1137.5  # H88: 0 -> 1137.5
$ws.Cells.Item(88, 9).Value = 600  # I88: 0 -> 600
$ws.Cells.Item(88, 10).Value = 1316.6666  # J88: 0 -> 1316.6666
$ws.Cells.Item(88, 11).Value = 600  # K88: 0 -> 600
$ws.Cells.Item(88, 12).Value = 1316.6666  # L88: 0 -> 1316.6666
$ws.Cells.Item(88, 13).Value = -194  # M88: None -> -194
$ws.Cells.Item(88, 14).Value = -2128.6666  # N88: None -> -2128.6666

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(91, 8).Value = 1137.5  # H91: 0 -> 1137.5
$ws.Cells.Item(91, 9).Value = 600  # I91: 0 -> 600
$ws.Cells.Item(91, 10).Value = 1316.6666  # J91: 0 -> 1316.6666
$ws.Cells.Item(91, 11).Value = 600  # K91: 0 -> 600
$ws.Cells.Item(91, 12).Value = 1316.6666  # L91: 0 -> 1316.6666
$ws.Cells.Item(91, 13).Value = 804  # M91: None -> 804
$ws.Cells.Item(91, 14).Value = -4124.6666  # N91: None -> -4124.6666

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 3073.3462  # H132: 3019.7407 -> 3073.3462
$ws.Cells.Item(132, 9).Value = 2356.1667  # I132: 2317.7368 -> 2356.1667
$ws.Cells.Item(132, 11).Value = 7068.500100000001  # K132: 6953.2104 -> 7068.500100000001
$ws.Cells.Item(132, 13).Value = -4538.500100000001  # M132: -4423.2104 -> -4538.500100000001

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 2749.2222  # H136: 2578.3845 -> 2749.2222
$ws.Cells.Item(136, 9).Value = 2593.125  # I136: 2392.818 -> 2593.125
$ws.Cells.Item(136, 10).Value = 3998  # J136: 3599 -> 3998
$ws.Cells.Item(136, 11).Value = 7779.375  # K136: 7178.454000000001 -> 7779.375
$ws.Cells.Item(136, 12).Value = 11994  # L136: 10797 -> 11994
$ws.Cells.Item(136, 13).Value = -5229.375  # M136: -4628.454000000001 -> -5229.375
$ws.Cells.Item(136, 14).Value = -17094  # N136: -15897 -> -17094

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2999.6667  # H134: 2997.1428 -> 2999.6667
$ws.Cells.Item(134, 9).Value = 2999.6667  # I134: 2997.1428 -> 2999.6667
$ws.Cells.Item(134, 11).Value = 8999.000100000001  # K134: 8991.428400000001 -> 8999.000100000001
$ws.Cells.Item(134, 13).Value = -6464.000100000001  # M134: -6456.428400000001 -> -6464.000100000001

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 4450  # H16: 0 -> 4450
$ws.Cells.Item(16, 9).Value = 900  # I16: 0 -> 900
$ws.Cells.Item(16, 10).Value = 8000  # J16: 0 -> 8000
$ws.Cells.Item(16, 11).Value = 900  # K16: 0 -> 900
$ws.Cells.Item(16, 12).Value = 8000  # L16: 0 -> 8000
$ws.Cells.Item(16, 13).Value = -613  # M16: None -> -613
$ws.Cells.Item(16, 14).Value = -8574  # N16: None -> -8574

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2940.889  # H31: 4349.0586 -> 2940.889
$ws.Cells.Item(31, 9).Value = 2541.9285  # I31: 2622.1538 -> 2541.9285
$ws.Cells.Item(31, 10).Value = 4337.25  # J31: 9961.5 -> 4337.25
$ws.Cells.Item(31, 11).Value = 2541.9285  # K31: 2622.1538 -> 2541.9285
$ws.Cells.Item(31, 12).Value = 4337.25  # L31: 9961.5 -> 4337.25
$ws.Cells.Item(31, 13).Value = -2246.9285  # M31: -2327.1538 -> -2246.9285
$ws.Cells.Item(31, 14).Value = -4927.25  # N31: -10551.5 -> -4927.25

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2940.889  # H34: 4349.0586 -> 2940.889
$ws.Cells.Item(34, 9).Value = 2541.9285  # I34: 2622.1538 -> 2541.9285
$ws.Cells.Item(34, 10).Value = 4337.25  # J34: 9961.5 -> 4337.25
$ws.Cells.Item(34, 11).Value = 2541.9285  # K34: 2622.1538 -> 2541.9285
$ws.Cells.Item(34, 12).Value = 4337.25  # L34: 9961.5 -> 4337.25
$ws.Cells.Item(34, 13).Value = -2339.9285  # M34: -2420.1538 -> -2339.9285
$ws.Cells.Item(34, 14).Value = -4741.25  # N34: -10365.5 -> -4741.25

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2142.923  # H58: 2129.1333 -> 2142.923
$ws.Cells.Item(58, 9).Value = 2142.923  # I58: 2129.1333 -> 2142.923
$ws.Cells.Item(58, 11).Value = 2142.923  # K58: 2129.1333 -> 2142.923
$ws.Cells.Item(58, 13).Value = -1939.923  # M58: -1926.1333 -> -1939.923

# CRP row 74
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(74, 8).Value = 50712.25  # H74: 56189.668 -> 50712.25
$ws.Cells.Item(74, 10).Value = 50712.25  # J74: 56189.668 -> 50712.25
$ws.Cells.Item(74, 12).Value = 50712.25  # L74: 56189.668 -> 50712.25
$ws.Cells.Item(74, 14).Value = -52460.25  # N74: -57937.668 -> -52460.25

# CRP row 77
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(77, 8).Value = 50712.25  # H77: 56189.668 -> 50712.25
$ws.Cells.Item(77, 10).Value = 50712.25  # J77: 56189.668 -> 50712.25
$ws.Cells.Item(77, 12).Value = 152136.75  # L77: 168569.004 -> 152136.75
$ws.Cells.Item(77, 14).Value = -160872.75  # N77: -177305.004 -> -160872.75

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 1750  # H105: 0 -> 1750
$ws.Cells.Item(105, 10).Value = 1750  # J105: 0 -> 1750
$ws.Cells.Item(105, 12).Value = 1750  # L105: 0 -> 1750
$ws.Cells.Item(105, 14).Value = -5244  # N105: None -> -5244

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 4450  # H113: 0 -> 4450
$ws.Cells.Item(113, 9).Value = 900  # I113: 0 -> 900
$ws.Cells.Item(113, 10).Value = 8000  # J113: 0 -> 8000
$ws.Cells.Item(113, 11).Value = 900  # K113: 0 -> 900
$ws.Cells.Item(113, 12).Value = 8000  # L113: 0 -> 8000
$ws.Cells.Item(113, 13).Value = 1270  # M113: None -> 1270
$ws.Cells.Item(113, 14).Value = -12340  # N113: None -> -12340

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 4484.6665  # H134: 5534.364 -> 4484.6665
$ws.Cells.Item(134, 9).Value = 4484.6665  # I134: 5534.364 -> 4484.6665
$ws.Cells.Item(134, 11).Value = 13453.9995  # K134: 16603.092 -> 13453.9995
$ws.Cells.Item(134, 13).Value = -10918.9995  # M134: -14068.092 -> -10918.9995

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 2142.923  # H136: 2129.1333 -> 2142.923
$ws.Cells.Item(136, 9).Value = 2142.923  # I136: 2129.1333 -> 2142.923
$ws.Cells.Item(136, 11).Value = 6428.768999999999  # K136: 6387.3999 -> 6428.768999999999
$ws.Cells.Item(136, 13).Value = -3878.768999999999  # M136: -3837.3999 -> -3878.768999999999

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 779.46155  # H5: 800.625 -> 779.46155
$ws.Cells.Item(5, 9).Value = 686.125  # I5: 687.25 -> 686.125
$ws.Cells.Item(5, 10).Value = 928.8  # J5: 914 -> 928.8
$ws.Cells.Item(5, 11).Value = 2058.375  # K5: 2061.75 -> 2058.375
$ws.Cells.Item(5, 12).Value = 2786.4  # L5: 2742 -> 2786.4
$ws.Cells.Item(5, 13).Value = -1946.375  # M5: -1949.75 -> -1946.375
$ws.Cells.Item(5, 14).Value = -3010.4  # N5: -2966 -> -3010.4

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 4178.625  # H132: 3858.75 -> 4178.625
$ws.Cells.Item(132, 9).Value = 975  # I132: 966 -> 975
$ws.Cells.Item(132, 10).Value = 5246.5  # J132: 5594.4 -> 5246.5
$ws.Cells.Item(132, 11).Value = 8775  # K132: 8694 -> 8775
$ws.Cells.Item(132, 12).Value = 47218.5  # L132: 50349.6 -> 47218.5
$ws.Cells.Item(132, 13).Value = -6245  # M132: -6164 -> -6245
$ws.Cells.Item(132, 14).Value = -52278.5  # N132: -55409.6 -> -52278.5

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 779.46155  # H135: 800.625 -> 779.46155
$ws.Cells.Item(135, 9).Value = 686.125  # I135: 687.25 -> 686.125
$ws.Cells.Item(135, 10).Value = 928.8  # J135: 914 -> 928.8
$ws.Cells.Item(135, 11).Value = 6175.125  # K135: 6185.25 -> 6175.125
$ws.Cells.Item(135, 12).Value = 8359.199999999999  # L135: 8226 -> 8359.199999999999
$ws.Cells.Item(135, 13).Value = -3640.125  # M135: -3650.25 -> -3640.125
$ws.Cells.Item(135, 14).Value = -13429.2  # N135: -13296 -> -13429.2

# GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 9).Value = 74.5  # I2: 68.59999999999999 -> 74.5
$ws.Cells.Item(2, 10).Value = 45  # J2: 0 -> 45
$ws.Cells.Item(2, 11).Value = 74.5  # K2: 68.59999999999999 -> 74.5
$ws.Cells.Item(2, 12).Value = 45  # L2: 0 -> 45
$ws.Cells.Item(2, 13).Value = 38.5  # M2: 44.40000000000001 -> 38.5
$ws.Cells.Item(2, 14).Value = -271  # N2: None -> -271

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 2065.6667  # H113: 2000 -> 2065.6667
$ws.Cells.Item(113, 9).Value = 1998.5  # I113: 2000 -> 1998.5
$ws.Cells.Item(113, 10).Value = 2200  # J113: 0 -> 2200
$ws.Cells.Item(113, 11).Value = 1998.5  # K113: 2000 -> 1998.5
$ws.Cells.Item(113, 12).Value = 2200  # L113: 0 -> 2200
$ws.Cells.Item(113, 13).Value = 171.5  # M113: 170 -> 171.5
$ws.Cells.Item(113, 14).Value = -6540  # N113: None -> -6540

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1097.25  # H122: 1141.6666 -> 1097.25
$ws.Cells.Item(122, 10).Value = 1498.5  # J122: 1498 -> 1498.5
$ws.Cells.Item(122, 12).Value = 4495.5  # L122: 4494 -> 4495.5
$ws.Cells.Item(122, 14).Value = -9395.5  # N122: -9394 -> -9395.5

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3727.2222  # H132: 3559.5 -> 3727.2222
$ws.Cells.Item(132, 9).Value = 3292.4285  # I132: 3206.6 -> 3292.4285
$ws.Cells.Item(132, 10).Value = 5249  # J132: 4618.2 -> 5249
$ws.Cells.Item(132, 11).Value = 9877.2855  # K132: 9619.799999999999 -> 9877.2855
$ws.Cells.Item(132, 12).Value = 15747  # L132: 13854.6 -> 15747
$ws.Cells.Item(132, 13).Value = -7347.2855  # M132: -7089.799999999999 -> -7347.2855
$ws.Cells.Item(132, 14).Value = -20807  # N132: -18914.6 -> -20807

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3899.6  # H7: 4249.5 -> 3899.6
$ws.Cells.Item(7, 9).Value = 3499.6667  # I7: 3999.5 -> 3499.6667
$ws.Cells.Item(7, 11).Value = 3499.6667  # K7: 3999.5 -> 3499.6667
$ws.Cells.Item(7, 13).Value = -3387.6667  # M7: -3887.5 -> -3387.6667

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3924.8667  # H22: 3931.3447 -> 3924.8667
$ws.Cells.Item(22, 9).Value = 3885.0588  # I22: 3894.3125 -> 3885.0588
$ws.Cells.Item(22, 11).Value = 3885.0588  # K22: 3894.3125 -> 3885.0588
$ws.Cells.Item(22, 13).Value = -3590.0588  # M22: -3599.3125 -> -3590.0588

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 3924.8667  # H27: 3931.3447 -> 3924.8667
$ws.Cells.Item(27, 9).Value = 3885.0588  # I27: 3894.3125 -> 3885.0588
$ws.Cells.Item(27, 11).Value = 3885.0588  # K27: 3894.3125 -> 3885.0588
$ws.Cells.Item(27, 13).Value = -3778.0588  # M27: -3787.3125 -> -3778.0588

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4952.5  # H40: 4966.3335 -> 4952.5
$ws.Cells.Item(40, 9).Value = 5000  # I40: 4997 -> 5000
$ws.Cells.Item(40, 11).Value = 5000  # K40: 4997 -> 5000
$ws.Cells.Item(40, 13).Value = -4864  # M40: -4861 -> -4864

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 3231.1667  # H61: 3777.4 -> 3231.1667
$ws.Cells.Item(61, 9).Value = 2877.8  # I61: 3472.25 -> 2877.8
$ws.Cells.Item(61, 11).Value = 2877.8  # K61: 3472.25 -> 2877.8
$ws.Cells.Item(61, 13).Value = -2675.8  # M61: -3270.25 -> -2675.8

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 3231.1667  # H113: 3777.4 -> 3231.1667
$ws.Cells.Item(113, 9).Value = 2877.8  # I113: 3472.25 -> 2877.8
$ws.Cells.Item(113, 11).Value = 2877.8  # K113: 3472.25 -> 2877.8
$ws.Cells.Item(113, 13).Value = -707.8000000000002  # M113: -1302.25 -> -707.8000000000002

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 1812.1428  # H122: 1841.1428 -> 1812.1428
$ws.Cells.Item(122, 9).Value = 1812.1428  # I122: 1848.6666 -> 1812.1428
$ws.Cells.Item(122, 10).Value = 0  # J122: 1796 -> 0
$ws.Cells.Item(122, 11).Value = 5436.428400000001  # K122: 5545.9998 -> 5436.428400000001
$ws.Cells.Item(122, 12).Value = 0  # L122: 5388 -> 0
$ws.Cells.Item(122, 13).Value = -2986.428400000001  # M122: -3095.9998 -> -2986.428400000001
$ws.Cells.Item(122, 14).ClearContents()  # N122: -10288 -> (removed)

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 3899.6  # H126: 4249.5 -> 3899.6
$ws.Cells.Item(126, 9).Value = 3499.6667  # I126: 3999.5 -> 3499.6667
$ws.Cells.Item(126, 11).Value = 10499.0001  # K126: 11998.5 -> 10499.0001
$ws.Cells.Item(126, 13).Value = -8029.000100000001  # M126: -9528.5 -> -8029.000100000001

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 33560.07  # H136: 35842.77 -> 33560.07
$ws.Cells.Item(136, 9).Value = 11641.5  # I136: 13192.8 -> 11641.5
$ws.Cells.Item(136, 11).Value = 34924.5  # K136: 39578.39999999999 -> 34924.5
$ws.Cells.Item(136, 13).Value = -32374.5  # M136: -37028.39999999999 -> -32374.5

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 898.125  # H81: 886.625 -> 898.125
$ws.Cells.Item(81, 9).Value = 898.125  # I81: 956.1429000000001 -> 898.125
$ws.Cells.Item(81, 10).Value = 0  # J81: 400 -> 0
$ws.Cells.Item(81, 11).Value = 1796.25  # K81: 1912.2858 -> 1796.25
$ws.Cells.Item(81, 12).Value = 0  # L81: 800 -> 0
$ws.Cells.Item(81, 13).Value = -735.25  # M81: -851.2858000000001 -> -735.25
$ws.Cells.Item(81, 14).ClearContents()  # N81: -2922 -> (removed)

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(84, 8).Value = 898.125  # H84: 886.625 -> 898.125
$ws.Cells.Item(84, 9).Value = 898.125  # I84: 956.1429000000001 -> 898.125
$ws.Cells.Item(84, 10).Value = 0  # J84: 400 -> 0
$ws.Cells.Item(84, 11).Value = 8981.25  # K84: 9561.429 -> 8981.25
$ws.Cells.Item(84, 12).Value = 0  # L84: 4000 -> 0
$ws.Cells.Item(84, 13).Value = -3677.25  # M84: -4257.429 -> -3677.25
$ws.Cells.Item(84, 14).ClearContents()  # N84: -14608 -> (removed)

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 726  # H113: 726.8570999999999 -> 726
$ws.Cells.Item(113, 9).Value = 726  # I113: 726.8570999999999 -> 726
$ws.Cells.Item(113, 11).Value = 2178  # K113: 2180.5713 -> 2178
$ws.Cells.Item(113, 13).Value = -8  # M113: -10.57129999999961 -> -8

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2696  # H122: 2795 -> 2696
$ws.Cells.Item(122, 9).Value = 2478.8333  # I122: 2494 -> 2478.8333
$ws.Cells.Item(122, 11).Value = 7436.499899999999  # K122: 7482 -> 7436.499899999999
$ws.Cells.Item(122, 13).Value = -4986.499899999999  # M122: -5032 -> -4986.499899999999

# WVR row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(135, 8).Value = 50000  # H135: 0 -> 50000
$ws.Cells.Item(135, 10).Value = 50000  # J135: 0 -> 50000
$ws.Cells.Item(135, 12).Value = 50000  # L135: 0 -> 50000
$ws.Cells.Item(135, 14).Value = -60140  # N135: None -> -60140

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 7975.421  # H136: 7980.2104 -> 7975.421
$ws.Cells.Item(136, 9).Value = 7975.421  # I136: 7980.2104 -> 7975.421
$ws.Cells.Item(136, 11).Value = 23926.263  # K136: 23940.6312 -> 23926.263
$ws.Cells.Item(136, 13).Value = -21376.263  # M136: -21390.6312 -> -21376.263
